$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 400
$ws.Range("C2").Value = 400

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 200
$ws.Range("C3").Value = 200

$ws.Range("A1:C1").Font.Bold = $true

$ws.Range("A4").Select()
